# Changing Dataschema in P1 (and P2) to ISCED-2011 Standard
# The "Categories" sheet lists allowed category values for each coded
# variable. EDU_LEVEL previously had 6 free-text categories; replace them
# with the 10 standard ISCED-2011 education levels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categories")

# The old EDU_LEVEL block occupies rows 4-9 (6 rows). The new block needs
# 10 rows, so insert 4 blank rows right after the existing block (i.e.
# before the following SMOKE_ST block which currently starts at row 10).
$ws.Rows("10:13").Insert()

# New ISCED-2011 education level categories (rows 4-13)
$levels = @(
    "Early Childhood Education",
    "Primary Education",
    "Lower Secondary Education",
    "Upper Secondary Education",
    "Post-secondary non-tertiary education",
    "Short-Cycle Tertiary Education",
    "Bachelor's or equivalent level",
    "Master's or equivalent level",
    "Doctoral or equivalent level",
    "Other"
)

for ($i = 0; $i -lt $levels.Length; $i++) {
    $row = 4 + $i
    $ws.Range("A$row").Value = "EDU_LEVEL"
    $ws.Range("B$row").Value = $levels[$i]
    $ws.Range("C$row").Value = $i
}

# Make "Categories" the active sheet/tab, with the new EDU_LEVEL block
# selected, mirroring the author's saved view state.
$ws.Activate()
$ws.Range("A4:C13").Select()
